# Scheduled market-data refresh: update the price/profit columns (H-N)
# for the Leve rows whose underlying item market prices moved, across the
# crafting-profession sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
#   H currentAveragePrice    I currentAveragePriceNQ  J currentAveragePriceHQ
#   K LevePriceNQ            L LevePriceHQ
#   M LeveProfitNQ           N LeveProfitHQ

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 8).Value = 7813581  # H98: was 17858412 -> 7813581
$ws.Cells.Item(98, 9).Value = 7813581  # I98: was 17858412 -> 7813581
$ws.Cells.Item(98, 11).Value = 7813581  # K98: was 17858412 -> 7813581
$ws.Cells.Item(98, 13).Value = -7812083  # M98: was -17856914 -> -7812083

$ws.Cells.Item(122, 8).Value = 7813581  # H122: was 17858412 -> 7813581
$ws.Cells.Item(122, 9).Value = 7813581  # I122: was 17858412 -> 7813581
$ws.Cells.Item(122, 11).Value = 23440743  # K122: was 53575236 -> 23440743
$ws.Cells.Item(122, 13).Value = -23438293  # M122: was -53572786 -> -23438293

$ws.Cells.Item(137, 8).Value = 1021.8033  # H137: was 997.64703 -> 1021.8033
$ws.Cells.Item(137, 9).Value = 829.1836499999999  # I137: was 876.11365 -> 829.1836499999999
$ws.Cells.Item(137, 10).Value = 1808.3334  # J137: was 1128.0731 -> 1808.3334
$ws.Cells.Item(137, 11).Value = 2487.55095  # K137: was 2628.34095 -> 2487.55095
$ws.Cells.Item(137, 12).Value = 5425.0002  # L137: was 3384.2193 -> 5425.0002
$ws.Cells.Item(137, 13).Value = 62.44905000000017  # M137: was -78.34094999999979 -> 62.44905000000017
$ws.Cells.Item(137, 14).Value = -10525.0002  # N137: was -8484.219300000001 -> -10525.0002

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 1353944  # H32: was 1677362.9 -> 1353944
$ws.Cells.Item(32, 9).Value = 1588946.8  # I32: was 2022195.2 -> 1588946.8
$ws.Cells.Item(32, 10).Value = 88544.38  # J32: was 96881.336 -> 88544.38
$ws.Cells.Item(32, 11).Value = 1588946.8  # K32: was 2022195.2 -> 1588946.8
$ws.Cells.Item(32, 12).Value = 88544.38  # L32: was 96881.336 -> 88544.38
$ws.Cells.Item(32, 13).Value = -1588659.8  # M32: was -2021908.2 -> -1588659.8
$ws.Cells.Item(32, 14).Value = -89118.38  # N32: was -97455.336 -> -89118.38

$ws.Cells.Item(132, 8).Value = 6188127.5  # H132: was 7303924.5 -> 6188127.5
$ws.Cells.Item(132, 9).Value = 8707605  # I132: was 10884382 -> 8707605
$ws.Cells.Item(132, 10).Value = 462042.53  # J132: was 484006.47 -> 462042.53
$ws.Cells.Item(132, 11).Value = 26122815  # K132: was 32653146 -> 26122815
$ws.Cells.Item(132, 12).Value = 1386127.59  # L132: was 1452019.41 -> 1386127.59
$ws.Cells.Item(132, 13).Value = -26120285  # M132: was -32650616 -> -26120285
$ws.Cells.Item(132, 14).Value = -1391187.59  # N132: was -1457079.41 -> -1391187.59

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 13186028  # H134: was 14737419 -> 13186028
$ws.Cells.Item(134, 9).Value = 19231936  # I134: was 20834776 -> 19231936
$ws.Cells.Item(134, 10).Value = 86562.164  # J134: was 103764.4 -> 86562.164
$ws.Cells.Item(134, 11).Value = 57695808  # K134: was 62504328 -> 57695808
$ws.Cells.Item(134, 12).Value = 259686.492  # L134: was 311293.2 -> 259686.492
$ws.Cells.Item(134, 13).Value = -57693273  # M134: was -62501793 -> -57693273
$ws.Cells.Item(134, 14).Value = -264756.492  # N134: was -316363.2 -> -264756.492

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 4058.8481  # H31: was 4247.16 -> 4058.8481
$ws.Cells.Item(31, 9).Value = 3337.2683  # I31: was 3640.973 -> 3337.2683
$ws.Cells.Item(31, 11).Value = 3337.2683  # K31: was 3640.973 -> 3337.2683
$ws.Cells.Item(31, 13).Value = -3042.2683  # M31: was -3345.973 -> -3042.2683

$ws.Cells.Item(34, 8).Value = 4058.8481  # H34: was 4247.16 -> 4058.8481
$ws.Cells.Item(34, 9).Value = 3337.2683  # I34: was 3640.973 -> 3337.2683
$ws.Cells.Item(34, 11).Value = 3337.2683  # K34: was 3640.973 -> 3337.2683
$ws.Cells.Item(34, 13).Value = -3135.2683  # M34: was -3438.973 -> -3135.2683

$ws.Cells.Item(36, 8).Value = 54615.9  # H36: was 43810.6 -> 54615.9
$ws.Cells.Item(36, 9).Value = 48000  # I36: was 45666.668 -> 48000
$ws.Cells.Item(36, 10).Value = 70053  # J36: was 41026.5 -> 70053
$ws.Cells.Item(36, 11).Value = 48000  # K36: was 45666.668 -> 48000
$ws.Cells.Item(36, 12).Value = 70053  # L36: was 41026.5 -> 70053
$ws.Cells.Item(36, 13).Value = -47612  # M36: was -45278.668 -> -47612
$ws.Cells.Item(36, 14).Value = -70829  # N36: was -41802.5 -> -70829

$ws.Cells.Item(40, 8).Value = 54615.9  # H40: was 43810.6 -> 54615.9
$ws.Cells.Item(40, 9).Value = 48000  # I40: was 45666.668 -> 48000
$ws.Cells.Item(40, 10).Value = 70053  # J40: was 41026.5 -> 70053
$ws.Cells.Item(40, 11).Value = 48000  # K40: was 45666.668 -> 48000
$ws.Cells.Item(40, 12).Value = 70053  # L40: was 41026.5 -> 70053
$ws.Cells.Item(40, 13).Value = -47840  # M40: was -45506.668 -> -47840
$ws.Cells.Item(40, 14).Value = -70373  # N40: was -41346.5 -> -70373

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(16, 8).Value = 801  # H16: was 2200 -> 801
$ws.Cells.Item(16, 9).Value = 600  # I16: was 800 -> 600
$ws.Cells.Item(16, 10).Value = 1002  # J16: was 5000 -> 1002
$ws.Cells.Item(16, 11).Value = 1800  # K16: was 2400 -> 1800
$ws.Cells.Item(16, 12).Value = 3006  # L16: was 15000 -> 3006
$ws.Cells.Item(16, 13).Value = -1627  # M16: was -2227 -> -1627
$ws.Cells.Item(16, 14).Value = -3352  # N16: was -15346 -> -3352

$ws.Cells.Item(20, 8).Value = 1730.2858  # H20: was 2168.3333 -> 1730.2858
$ws.Cells.Item(20, 9).Value = 100  # I20: was 0 -> 100
$ws.Cells.Item(20, 10).Value = 2002  # J20: was 2168.3333 -> 2002
$ws.Cells.Item(20, 11).Value = 300  # K20: was 0 -> 300
$ws.Cells.Item(20, 12).Value = 6006  # L20: was 6504.999899999999 -> 6006
$ws.Cells.Item(20, 13).Value = -73  # M20: was (blank) -> -73
$ws.Cells.Item(20, 14).Value = -6460  # N20: was -6958.999899999999 -> -6460

$ws.Cells.Item(21, 8).Value = 699.3333  # H21: was 0 -> 699.3333
$ws.Cells.Item(21, 9).Value = 699.3333  # I21: was 0 -> 699.3333
$ws.Cells.Item(21, 11).Value = 2097.9999  # K21: was 0 -> 2097.9999
$ws.Cells.Item(21, 13).Value = -1924.9999  # M21: was (blank) -> -1924.9999

$ws.Cells.Item(22, 8).Value = 1465  # H22: was 1107.8462 -> 1465
$ws.Cells.Item(22, 9).Value = 816  # I22: was 933.3333 -> 816
$ws.Cells.Item(22, 10).Value = 1928.5714  # J22: was 1257.4286 -> 1928.5714
$ws.Cells.Item(22, 11).Value = 2448  # K22: was 2799.9999 -> 2448
$ws.Cells.Item(22, 12).Value = 5785.7142  # L22: was 3772.2858 -> 5785.7142
$ws.Cells.Item(22, 13).Value = -2279  # M22: was -2630.9999 -> -2279
$ws.Cells.Item(22, 14).Value = -6123.7142  # N22: was -4110.2858 -> -6123.7142

$ws.Cells.Item(27, 8).Value = 1465  # H27: was 1107.8462 -> 1465
$ws.Cells.Item(27, 9).Value = 816  # I27: was 933.3333 -> 816
$ws.Cells.Item(27, 10).Value = 1928.5714  # J27: was 1257.4286 -> 1928.5714
$ws.Cells.Item(27, 11).Value = 2448  # K27: was 2799.9999 -> 2448
$ws.Cells.Item(27, 12).Value = 5785.7142  # L27: was 3772.2858 -> 5785.7142
$ws.Cells.Item(27, 13).Value = -2346  # M27: was -2697.9999 -> -2346
$ws.Cells.Item(27, 14).Value = -5989.7142  # N27: was -3976.2858 -> -5989.7142

$ws.Cells.Item(58, 8).Value = 5496504.5  # H58: was 7694177 -> 5496504.5
$ws.Cells.Item(58, 10).Value = 6412422  # J58: was 9617472 -> 6412422
$ws.Cells.Item(58, 12).Value = 19237266  # L58: was 28852416 -> 19237266
$ws.Cells.Item(58, 14).Value = -19237522  # N58: was -28852672 -> -19237522

$ws.Cells.Item(61, 8).Value = 100  # H61: was 800 -> 100
$ws.Cells.Item(61, 10).Value = 0  # J61: was 1150 -> 0
$ws.Cells.Item(61, 12).Value = 0  # L61: was 3450 -> 0
$ws.Cells.Item(61, 14).ClearContents()  # N61: was -3880 -> (cleared)

$ws.Cells.Item(97, 8).Value = 3454  # H97: was 5438.75 -> 3454
$ws.Cells.Item(97, 9).Value = 950  # I97: was 0 -> 950
$ws.Cells.Item(97, 10).Value = 3632.8572  # J97: was 5438.75 -> 3632.8572
$ws.Cells.Item(97, 11).Value = 2850  # K97: was 0 -> 2850
$ws.Cells.Item(97, 12).Value = 10898.5716  # L97: was 16316.25 -> 10898.5716
$ws.Cells.Item(97, 13).Value = -2354  # M97: was (blank) -> -2354
$ws.Cells.Item(97, 14).Value = -11890.5716  # N97: was -17308.25 -> -11890.5716

$ws.Cells.Item(131, 8).Value = 873.4865  # H131: was 902.6857 -> 873.4865
$ws.Cells.Item(131, 9).Value = 513.4286  # I131: was 563.1667 -> 513.4286
$ws.Cells.Item(131, 10).Value = 1092.6522  # J131: was 1079.826 -> 1092.6522
$ws.Cells.Item(131, 11).Value = 1540.2858  # K131: was 1689.5001 -> 1540.2858
$ws.Cells.Item(131, 12).Value = 3277.9566  # L131: was 3239.478 -> 3277.9566
$ws.Cells.Item(131, 13).Value = 3499.7142  # M131: was 3350.4999 -> 3499.7142
$ws.Cells.Item(131, 14).Value = -13357.9566  # N131: was -13319.478 -> -13357.9566

$ws.Cells.Item(140, 8).Value = 5089.2104  # H140: was 6114.4443 -> 5089.2104
$ws.Cells.Item(140, 9).Value = 1875.6  # I140: was 2147.1428 -> 1875.6
$ws.Cells.Item(140, 10).Value = 11269.23  # J140: was 20000 -> 11269.23
$ws.Cells.Item(140, 11).Value = 5626.799999999999  # K140: was 6441.428400000001 -> 5626.799999999999
$ws.Cells.Item(140, 12).Value = 33807.69  # L140: was 60000 -> 33807.69
$ws.Cells.Item(140, 13).Value = -446.7999999999993  # M140: was -1261.428400000001 -> -446.7999999999993
$ws.Cells.Item(140, 14).Value = -44167.69  # N140: was -70360 -> -44167.69

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(95, 8).Value = 9477.714  # H95: was 3923.3333 -> 9477.714
$ws.Cells.Item(95, 10).Value = 9477.714  # J95: was 3923.3333 -> 9477.714
$ws.Cells.Item(95, 12).Value = 9477.714  # L95: was 3923.3333 -> 9477.714
$ws.Cells.Item(95, 14).Value = -14969.714  # N95: was -9415.3333 -> -14969.714

$ws.Cells.Item(132, 8).Value = 18389.58  # H132: was 20640.254 -> 18389.58
$ws.Cells.Item(132, 9).Value = 1516.0834  # I132: was 1687.9667 -> 1516.0834
$ws.Cells.Item(132, 10).Value = 41752.883  # J132: was 43383 -> 41752.883
$ws.Cells.Item(132, 11).Value = 4548.2502  # K132: was 5063.9001 -> 4548.2502
$ws.Cells.Item(132, 12).Value = 125258.649  # L132: was 130149 -> 125258.649
$ws.Cells.Item(132, 13).Value = -2018.2502  # M132: was -2533.9001 -> -2018.2502
$ws.Cells.Item(132, 14).Value = -130318.649  # N132: was -135209 -> -130318.649

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(68, 8).Value = 2006.5883  # H68: was 1953.5714 -> 2006.5883
$ws.Cells.Item(68, 9).Value = 1749.3334  # I68: was 1730.9231 -> 1749.3334
$ws.Cells.Item(68, 10).Value = 2624  # J68: was 2315.375 -> 2624
$ws.Cells.Item(68, 11).Value = 1749.3334  # K68: was 1730.9231 -> 1749.3334
$ws.Cells.Item(68, 12).Value = 2624  # L68: was 2315.375 -> 2624
$ws.Cells.Item(68, 13).Value = -1000.3334  # M68: was -981.9231 -> -1000.3334
$ws.Cells.Item(68, 14).Value = -4122  # N68: was -3813.375 -> -4122

$ws.Cells.Item(71, 8).Value = 2006.5883  # H71: was 1953.5714 -> 2006.5883
$ws.Cells.Item(71, 9).Value = 1749.3334  # I71: was 1730.9231 -> 1749.3334
$ws.Cells.Item(71, 10).Value = 2624  # J71: was 2315.375 -> 2624
$ws.Cells.Item(71, 11).Value = 8746.666999999999  # K71: was 8654.6155 -> 8746.666999999999
$ws.Cells.Item(71, 12).Value = 13120  # L71: was 11576.875 -> 13120
$ws.Cells.Item(71, 13).Value = -5002.666999999999  # M71: was -4910.6155 -> -5002.666999999999
$ws.Cells.Item(71, 14).Value = -20608  # N71: was -19064.875 -> -20608

$ws.Cells.Item(132, 8).Value = 224062.61  # H132: was 281152.66 -> 224062.61
$ws.Cells.Item(132, 9).Value = 56716.805  # I132: was 70265.07000000001 -> 56716.805
$ws.Cells.Item(132, 10).Value = 558754.25  # J132: was 717991.2 -> 558754.25
$ws.Cells.Item(132, 11).Value = 170150.415  # K132: was 210795.21 -> 170150.415
$ws.Cells.Item(132, 12).Value = 1676262.75  # L132: was 2153973.6 -> 1676262.75
$ws.Cells.Item(132, 13).Value = -167620.415  # M132: was -208265.21 -> -167620.415
$ws.Cells.Item(132, 14).Value = -1681322.75  # N132: was -2159033.6 -> -1681322.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 3704.5476  # H132: was 4668.6177 -> 3704.5476
$ws.Cells.Item(132, 9).Value = 932.14813  # I132: was 1157.0834 -> 932.14813
$ws.Cells.Item(132, 10).Value = 8694.866  # J132: was 13096.3 -> 8694.866
$ws.Cells.Item(132, 11).Value = 2796.44439  # K132: was 3471.2502 -> 2796.44439
$ws.Cells.Item(132, 12).Value = 26084.598  # L132: was 39288.89999999999 -> 26084.598
$ws.Cells.Item(132, 13).Value = -266.4443900000001  # M132: was -941.2501999999999 -> -266.4443900000001
$ws.Cells.Item(132, 14).Value = -31144.598  # N132: was -44348.89999999999 -> -31144.598

$ws.Cells.Item(136, 8).Value = 1733679  # H136: was 2469121.2 -> 1733679
$ws.Cells.Item(136, 9).Value = 2382227.8  # I136: was 3970262.5 -> 2382227.8
$ws.Cells.Item(136, 10).Value = 589181.2  # J136: was 667751.7 -> 589181.2
$ws.Cells.Item(136, 11).Value = 7146683.399999999  # K136: was 11910787.5 -> 7146683.399999999
$ws.Cells.Item(136, 12).Value = 1767543.6  # L136: was 2003255.1 -> 1767543.6
$ws.Cells.Item(136, 13).Value = -7144133.399999999  # M136: was -11908237.5 -> -7144133.399999999
$ws.Cells.Item(136, 14).Value = -1772643.6  # N136: was -2008355.1 -> -1772643.6
